# Rename sheet "EDA" -> "EDAPrimer2" and refresh the measurement data on both
# the EDAPrimer and EDAPrimer2 (formerly EDA) sheets: the summary cells (B4/B5)
# get new median values, the existing per-run rows (12/13) get updated
# figures, and two additional per-run rows (14/15/16) are appended.

$wb = $excel.ActiveWorkbook

# --- Rename the EDA sheet -------------------------------------------------
$edaSheet = $wb.Worksheets.Item("EDA")
$edaSheet.Name = "EDAPrimer2"

$primerSheet = $wb.Worksheets.Item("EDAPrimer")
$primer2Sheet = $wb.Worksheets.Item("EDAPrimer2")

# The sheet's title cell (A1) mirrors the sheet name - keep it in sync too.
$primer2Sheet.Range("A1").Value = "EDAPrimer2"

# --- EDAPrimer sheet updates ----------------------------------------------
$primerSheet.Range("B4").Value = 162388.0
$primerSheet.Range("B5").Value = 19.671521035598705

$primerSheet.Range("B12").Value = 149127.0
$primerSheet.Range("C12").Value = 1009.0
$primerSheet.Range("E12").Value = 19.857605177993527

$primerSheet.Range("B13").Value = 143894.0
$primerSheet.Range("C13").Value = 761.0
$primerSheet.Range("E13").Value = 17.733009708737864

# New rows 14-16: copy formatting down from row 13, then fill in values.
$primerSheet.Range("A13:Q13").Copy() | Out-Null
$primerSheet.Range("A14:Q14").PasteSpecial(-4122) | Out-Null
$primerSheet.Range("A15:Q15").PasteSpecial(-4122) | Out-Null
$primerSheet.Range("A16:Q16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$primerRows = @{
    14 = @(2.0, 162388.0, 1088.0, 0.0, 19.218446601941746)
    15 = @(3.0, 181338.0, 941.0, 0.0, 21.723300970873787)
    16 = @(4.0, 182415.0, 775.0, 0.0, 19.671521035598705)
}
foreach ($r in $primerRows.Keys) {
    $vals = $primerRows[$r]
    $primerSheet.Range("A$r").Value = $vals[0]
    $primerSheet.Range("B$r").Value = $vals[1]
    $primerSheet.Range("C$r").Value = $vals[2]
    $primerSheet.Range("D$r").Value = $vals[3]
    $primerSheet.Range("E$r").Value = $vals[4]
    for ($col = 6; $col -le 17; $col++) {
        $primerSheet.Cells.Item($r, $col).Value = 0.0
    }
}

# --- EDAPrimer2 (formerly EDA) sheet updates ------------------------------
$primer2Sheet.Range("B4").Value = 165773.0
$primer2Sheet.Range("B5").Value = 18.681229773462782

$primer2Sheet.Range("B12").Value = 165773.0
$primer2Sheet.Range("C12").Value = 1132.0
$primer2Sheet.Range("E12").Value = 18.741100323624597

$primer2Sheet.Range("B13").Value = 174375.0
$primer2Sheet.Range("C13").Value = 1399.0
$primer2Sheet.Range("E13").Value = 19.45145631067961

# New rows 14-16: copy formatting down from row 13, then fill in values.
$primer2Sheet.Range("A13:Q13").Copy() | Out-Null
$primer2Sheet.Range("A14:Q14").PasteSpecial(-4122) | Out-Null
$primer2Sheet.Range("A15:Q15").PasteSpecial(-4122) | Out-Null
$primer2Sheet.Range("A16:Q16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$primer2Rows = @{
    14 = @(2.0, 166449.0, 1119.0, 0.0, 18.681229773462782)
    15 = @(3.0, 133298.0, 668.0, 0.0, 15.21359223300971)
    16 = @(4.0, 140785.0, 715.0, 0.0, 17.401294498381876)
}
foreach ($r in $primer2Rows.Keys) {
    $vals = $primer2Rows[$r]
    $primer2Sheet.Range("A$r").Value = $vals[0]
    $primer2Sheet.Range("B$r").Value = $vals[1]
    $primer2Sheet.Range("C$r").Value = $vals[2]
    $primer2Sheet.Range("D$r").Value = $vals[3]
    $primer2Sheet.Range("E$r").Value = $vals[4]
    for ($col = 6; $col -le 17; $col++) {
        $primer2Sheet.Cells.Item($r, $col).Value = 0.0
    }
}

Write-Host "Measurements workbook updated."
